$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-07 05:47:39"
$ws.Range("N2").Value = "-3.4 °C 5:29 TU"
$ws.Range("O2").Value = "-1.5 °C"
$ws.Range("E3").Value = "2026-02-07 05:47:41"
$ws.Range("N3").Value = "-7.1 °C 5:07 TU"
$ws.Range("O3").Value = "-5.8 °C"
$ws.Range("E4").Value = "2026-02-07 05:47:44"
$ws.Range("J4").Value = "1001.1 hPa"
$ws.Range("O4").Value = "11.4 °C"
$ws.Range("E5").Value = "2026-02-07 05:47:46"
$ws.Range("H5").Value = "'73%"
$ws.Range("J5").Value = "1001.2 hPa"
$ws.Range("N5").Value = "6.6 °C 5:28 TU"
$ws.Range("O5").Value = "8.7 °C"
$ws.Range("E6").Value = "2026-02-07 05:47:49"
$ws.Range("K6").Value = "-0.1 MJ/m2"
$ws.Range("O6").Value = "11.8 °C"
$ws.Range("E7").Value = "2026-02-07 05:47:51"
$ws.Range("H7").Value = "'75%"
$ws.Range("J7").Value = "1002.6 hPa"
$ws.Range("E8").Value = "2026-02-07 05:47:53"
$ws.Range("K8").Value = "-0.1 MJ/m2"
$ws.Range("N8").Value = "2.4 °C 5:16 TU"
$ws.Range("O8").Value = "4.0 °C"
$ws.Range("E9").Value = "2026-02-07 05:47:56"
$ws.Range("N9").Value = "-0.4 °C 5:18 TU"
$ws.Range("O9").Value = "1.7 °C"
$ws.Range("E10").Value = "2026-02-07 05:47:58"
$ws.Range("M10").Value = "8.8 °C 5:13 TU"
$ws.Range("O10").Value = "7.2 °C"
$ws.Range("E11").Value = "2026-02-07 05:48:00"
$ws.Range("J11").Value = "1005.5 hPa"
$ws.Range("E12").Value = "2026-02-07 05:48:02"
$ws.Range("N12").Value = "7.0 °C 5:00 TU"
$ws.Range("O12").Value = "9.6 °C"
$ws.Range("E13").Value = "2026-02-07 05:48:05"
$ws.Range("H13").Value = "'89%"
$ws.Range("M13").Value = "10.7 °C 5:23 TU"
$ws.Range("O13").Value = "7.5 °C"
$ws.Range("E14").Value = "2026-02-07 05:48:07"
$ws.Range("H14").Value = "'75%"
$ws.Range("K14").Value = "-0.1 MJ/m2"
$ws.Range("E15").Value = "2026-02-07 05:48:09"
$ws.Range("H15").Value = "'86%"
$ws.Range("J15").Value = "1001.5 hPa"
$ws.Range("N15").Value = "3.4 °C 5:26 TU"
$ws.Range("O15").Value = "6.4 °C"
$ws.Range("E16").Value = "2026-02-07 05:48:12"
$ws.Range("H16").Value = "'91%"
$ws.Range("E17").Value = "2026-02-07 05:48:14"
$ws.Range("J17").Value = "1004.8 hPa"
$ws.Range("L17").Value = "9.0 km/h - 243º 5:28 TU"
$ws.Range("N17").Value = "2.5 °C 5:29 TU"
$ws.Range("O17").Value = "3.3 °C"
$ws.Range("E18").Value = "2026-02-07 05:48:17"
$ws.Range("H18").Value = "'93%"
$ws.Range("N18").Value = "-9.5 °C 5:29 TU"
$ws.Range("O18").Value = "-7.4 °C"
$ws.Range("E19").Value = "2026-02-07 05:48:19"
$ws.Range("J19").Value = "1006.0 hPa"
$ws.Range("N19").Value = "2.8 °C 5:29 TU"
$ws.Range("O19").Value = "4.4 °C"
$ws.Range("E20").Value = "2026-02-07 05:48:22"
$ws.Range("N20").Value = "-5.7 °C 5:19 TU"
$ws.Range("O20").Value = "-4.6 °C"
$ws.Range("E21").Value = "2026-02-07 05:48:24"
$ws.Range("H21").Value = "'75%"
$ws.Range("J21").Value = "1001.7 hPa"
$ws.Range("N21").Value = "2.6 °C 5:04 TU"
$ws.Range("O21").Value = "6.7 °C"
$ws.Range("E22").Value = "2026-02-07 05:48:26"
$ws.Range("E23").Value = "2026-02-07 05:48:29"
$ws.Range("H23").Value = "'98%"
$ws.Range("J23").Value = "1001.4 hPa"
$ws.Range("E24").Value = "2026-02-07 05:48:31"
$ws.Range("J24").Value = "1000.7 hPa"
$ws.Range("L24").Value = "62.6 km/h - 339º 5:29 TU"
$ws.Range("N24").Value = "9.8 °C 5:20 TU"
$ws.Range("E25").Value = "2026-02-07 05:48:33"
$ws.Range("H25").Value = "'98%"
$ws.Range("E26").Value = "2026-02-07 05:48:36"
$ws.Range("N26").Value = "-4.9 °C 5:21 TU"
$ws.Range("O26").Value = "-1.9 °C"
$ws.Range("E27").Value = "2026-02-07 05:48:38"
$ws.Range("H27").Value = "'95%"
$ws.Range("M27").Value = "10.1 °C 5:26 TU"
$ws.Range("O27").Value = "8.6 °C"
$ws.Range("E28").Value = "2026-02-07 05:48:41"
$ws.Range("H28").Value = "'88%"
$ws.Range("J28").Value = "1003.9 hPa"
$ws.Range("N28").Value = "1.7 °C 5:01 TU"
$ws.Range("O28").Value = "3.2 °C"
$ws.Range("E29").Value = "2026-02-07 05:48:43"
$ws.Range("K29").Value = "-0.1 MJ/m2"
$ws.Range("O29").Value = "11.0 °C"
$ws.Range("E30").Value = "2026-02-07 05:48:45"
$ws.Range("O30").Value = "-4.9 °C"
$ws.Range("E31").Value = "2026-02-07 05:48:48"
$ws.Range("J31").Value = "1005.8 hPa"
$ws.Range("N31").Value = "3.0 °C 5:29 TU"
$ws.Range("O31").Value = "3.6 °C"
$ws.Range("E32").Value = "2026-02-07 05:48:50"
$ws.Range("H32").Value = "'59%"
$ws.Range("J32").Value = "1004.3 hPa"
$ws.Range("E33").Value = "2026-02-07 05:48:53"
$ws.Range("N33").Value = "5.3 °C 5:25 TU"
$ws.Range("O33").Value = "7.0 °C"
$ws.Range("E34").Value = "2026-02-07 05:48:55"
$ws.Range("H34").Value = "'77%"
$ws.Range("N34").Value = "4.8 °C 5:10 TU"
$ws.Range("O34").Value = "6.3 °C"
$ws.Range("E35").Value = "2026-02-07 05:48:57"
$ws.Range("H35").Value = "'92%"
$ws.Range("O35").Value = "-5.7 °C"
$ws.Range("E36").Value = "2026-02-07 05:49:00"
$ws.Range("J36").Value = "1006.6 hPa"
$ws.Range("O36").Value = "4.6 °C"
